$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price (D) and Volume (E) columns stay as Text so that
# values like "314.48" or "1.000" are not reinterpreted as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "28.412.46"
$ws.Range("E2").Value = "  +4.12%  "
$ws.Range("D3").Value = "1.796.19"
$ws.Range("E3").Value = "  +1.15%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "314.48"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "0.5458"
$ws.Range("E7").Value = "  +6.15%  "
$ws.Range("D8").Value = "0.3831"
$ws.Range("E8").Value = "  +4.00%  "
$ws.Range("D9").Value = "0.07583"
$ws.Range("E9").Value = "  +2.65%  "
$ws.Range("D10").Value = "42.62"
$ws.Range("E10").Value = "  -0.41%  "
$ws.Range("E11").Value = "  +3.45%  "
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").Value = "21.15"
$ws.Range("E13").Value = "  +3.35%  "
$ws.Range("D14").Value = "6.185"
$ws.Range("E14").Value = "  +1.94%  "
$ws.Range("D15").Value = "7.406"
$ws.Range("E15").Value = "  +6.47%  "
$ws.Range("D16").Value = "1.795.52"
$ws.Range("E16").Value = "  +1.55%  "
$ws.Range("D17").Value = "91.50"
$ws.Range("E17").Value = "  +2.66%  "
$ws.Range("E18").Value = "  +2.47%  "
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("D20").Value = "1.000"
$ws.Range("E21").Value = "  +3.32%  "
$ws.Range("D22").Value = "5.960"
$ws.Range("E22").Value = "  +2.28%  "
$ws.Range("D23").Value = "28.409.90"
$ws.Range("E23").Value = "  +3.99%  "
$ws.Range("E24").Value = "  +2.04%  "
$ws.Range("D25").Value = "2.128"
$ws.Range("E25").Value = "  +0.46%  "
$ws.Range("D26").Value = "159.73"
$ws.Range("E26").Value = "  +3.27%  "
$ws.Range("D27").Value = "20.70"
$ws.Range("E27").Value = "  +2.57%  "
$ws.Range("E28").Value = "  +2.90%  "
$ws.Range("D29").Value = "2.002.35"
$ws.Range("E29").Value = "  +1.43%  "
$ws.Range("D30").Value = "123.34"
$ws.Range("E30").Value = "  +1.83%  "
$ws.Range("D31").Value = "1.127"
$ws.Range("E31").Value = "  +5.57%  "
$ws.Range("D32").Value = "0.1028"
$ws.Range("E32").Value = "  +5.40%  "
$ws.Range("D33").Value = "5.757"
$ws.Range("E33").Value = "  +3.33%  "
$ws.Range("D34").Value = "3.683"
$ws.Range("E34").Value = "  +1.53%  "
$ws.Range("D37").Value = "0.02324"
$ws.Range("E37").Value = "  +3.61%  "
$ws.Range("D38").Value = "5.166"
$ws.Range("E38").Value = "  +6.71%  "
$ws.Range("D39").Value = "8.772"
$ws.Range("E39").Value = "  +8.43%  "
$ws.Range("E40").Value = "  +3.85%  "
$ws.Range("D41").Value = "0.6398"
$ws.Range("E41").Value = "  +4.06%  "
$ws.Range("D42").Value = "0.9997"
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").Value = "1.410"
$ws.Range("E43").Value = "  -1.85%  "
$ws.Range("D44").Value = "1.160"
$ws.Range("E44").Value = "  +1.92%  "
$ws.Range("D45").Value = "13.52"
$ws.Range("E45").Value = "  +2.95%  "
$ws.Range("D46").Value = "0.5978"
$ws.Range("E46").Value = "  +3.74%  "
$ws.Range("D47").Value = "3.665"
$ws.Range("E47").Value = "  +0.97%  "
$ws.Range("D48").Value = "126.56"
$ws.Range("E48").Value = "  +4.37%  "
$ws.Range("D49").Value = "2.003"
$ws.Range("E49").Value = "  +6.26%  "
$ws.Range("D50").Value = "1.151"
$ws.Range("E50").Value = "  +3.40%  "
$ws.Range("D51").Value = "0.06954"
$ws.Range("E51").Value = "  +3.70%  "

# Row 35/36 swap: Algorand <-> Hedera
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "0.06793"
$ws.Range("E35").Value = "  +13.71%  "

$ws.Range("B36").Value = "Algorand"
$ws.Range("C36").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D36").Value = "0.2347"
$ws.Range("E36").Value = "  +16.12%  "
